$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 89 (shifts existing rows 89..109 down to 90..110)
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with this week's data point
$ws.Range("A89").Value = 10
$ws.Range("B89").Value = "Vega Modelo de Temuco"
$ws.Range("C89").Value = "La Araucanía"
$ws.Range("D89").Value = 44782
$ws.Range("E89").Value = 9
$ws.Range("F89").Value = 100112035
$ws.Range("G89").Value = "Bruselas (repollito)"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 30
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = 25000
$ws.Range("N89").Value = "$/malla 10 kilos"
$ws.Range("O89").Value = "Provincia de Quillota"
$ws.Range("P89").Value = 2500
$ws.Range("Q89").Value = 10
$ws.Range("R89").Value = "Hortaliza"
